$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: replace a table cell's visible text while preserving the run's
# formatting (rPr) and without touching the cell-end / paragraph-end
# marks. Cell.Range.Text includes a trailing paragraph mark + cell mark
# (2 characters), and when a cell holds more than one run, assigning
# straight to Cell.Range.Text only overwrites the first run and leaves
# any trailing runs behind -- so we explicitly size the replacement
# range to the cell's real text span (excluding the two end marks).
function Set-CellText($cell, $newText) {
    $rng = $cell.Range
    $len = $rng.Text.Length
    $textRng = $d.Range($rng.Start, $rng.Start + $len - 2)
    $textRng.Text = $newText
}

# The "Use case description" (URS) column and the "UI" (traceability)
# column for rows No=15..19 were reviewed/re-ordered, while the row
# numbers (No, URS.ID, UC, AD, SD) and the SRS.ID column stayed where
# they were. Net effect: each row's URS text + UI value is replaced
# with what used to belong to the next row down, and row 19 receives
# what used to be in row 15.

# Row No=15 (table row 16): "User can register..." -> "Registered user can login..."
$row = $t.Rows.Item(16)
Set-CellText $row.Cells.Item(3) "Registered user can login to the website."
Set-CellText $row.Cells.Item(8) "01"

# Row No=16 (table row 17): "Registered user can login..." -> "Registered user can log out..."
$row = $t.Rows.Item(17)
Set-CellText $row.Cells.Item(3) "Registered user can log out to the website."
Set-CellText $row.Cells.Item(8) "02,03,04,05,15,18,19"

# Row No=17 (table row 18): "Registered user can log out..." -> "Registered user can edit..."
$row = $t.Rows.Item(18)
Set-CellText $row.Cells.Item(3) "Registered user can edit their own registered information."
Set-CellText $row.Cells.Item(8) "18"

# Row No=18 (table row 19): "Registered user can edit..." -> "Registered user can get news..."
$row = $t.Rows.Item(19)
Set-CellText $row.Cells.Item(3) "Registered user can get news from admin by their registered e-mail."
Set-CellText $row.Cells.Item(8) "-"

# Row No=19 (table row 20): "Registered user can get news..." -> "User can register..."
$row = $t.Rows.Item(20)
Set-CellText $row.Cells.Item(3) "User can register to the website."
Set-CellText $row.Cells.Item(8) "17"

# Row No=29 (table row 30): SRS.ID gains a leading "65," reference.
$row = $t.Rows.Item(30)
Set-CellText $row.Cells.Item(4) "65,114,121,122"
